## Remove the SharePoint/document-library "custom XML" parts that were
## bundled into the package (customXml/item2.xml, item3.xml, item4.xml and
## their companion itemProps2-4.xml datastore items). These are the
## auto-generated "Document" content-type schema, the SharePoint
## DocumentLibraryForm template and the documentManagement properties
## part - none of them are referenced by the document body, so removing
## them leaves the visible content untouched.

$d = $word.ActiveDocument

# The three orphaned custom XML parts, identified by the well known
# itemID GUIDs recorded in their companion itemProps*.xml datastore
# items (and, as a fallback, by the root-element namespace URI of the
# part itself).
$idsToRemove = @(
    "{397287B1-721A-475D-8EFE-0CC58780EB97}",  # customXml/item2.xml  (ct:contentTypeSchema)
    "{A30FB459-7519-4D92-A060-DE50C1DAD0B1}",  # customXml/item3.xml  (FormTemplates)
    "{3CF42E69-F47B-43D5-B871-45F59531B940}"   # customXml/item4.xml  (p:properties)
)

$namespacesToRemove = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)

$parts = $d.CustomXMLParts

# First pass: try a direct, precise removal by itemID.
foreach ($id in $idsToRemove) {
    try {
        $part = $parts.SelectByID($id)
        if ($part -ne $null) {
            $part.Delete()
        }
    } catch {
    }
}

# Second pass: fall back to namespace-based lookup/removal in case
# SelectByID isn't available, walking any matching selection.
foreach ($ns in $namespacesToRemove) {
    try {
        $sel = $parts.SelectByNamespace($ns)
        if ($sel -ne $null) {
            for ($i = $sel.Count; $i -ge 1; $i--) {
                try { $sel.Item($i).Delete() } catch {}
            }
        }
    } catch {
    }
}

# Third pass: brute-force walk of the whole collection (descending, so
# deleting doesn't shift the indices we still need to visit), matching
# on either the itemID or the part's namespace URI.
try {
    for ($i = $parts.Count; $i -ge 1; $i--) {
        try {
            $p = $parts.Item($i)
            if ($p -eq $null) { continue }

            $pid = $null
            try { $pid = $p.Id } catch {}
            if ($pid -eq $null) { try { $pid = $p.ID } catch {} }

            $pns = $null
            try { $pns = $p.NamespaceURI } catch {}

            if (($idsToRemove -contains $pid) -or ($namespacesToRemove -contains $pns)) {
                $p.Delete()
            }
        } catch {
        }
    }
} catch {
}

Write-Host "Done removing orphaned custom XML parts."
